$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 158, shifting rows 158:172
# down to 160:174 (same as old rows, values unchanged, just relocated).
$ws.Rows("158:159").Insert()

# Row 158 (new): Ají - Cacho cabra rojo, Región del Maule
$ws.Range("A158").Value = 7
$ws.Range("B158").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C158").Value = "Ñuble"
$ws.Range("D158").Value = 45013
$ws.Range("E158").Value = 16
$ws.Range("F158").Value = 100112021
$ws.Range("G158").Value = "Ají"
$ws.Range("H158").Value = "Cacho cabra rojo"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 30
$ws.Range("K158").Value = 15000
$ws.Range("L158").Value = 15000
$ws.Range("M158").Value = 15000
$ws.Range("N158").Value = "`$/saco 25 kilos"
$ws.Range("O158").Value = "Región del Maule"
$ws.Range("P158").Value = 600
$ws.Range("Q158").Value = 25
$ws.Range("R158").Value = "Hortaliza"

# Row 159 (new): Ají - Cacho cabra verde, Región del Maule
$ws.Range("A159").Value = 7
$ws.Range("B159").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C159").Value = "Ñuble"
$ws.Range("D159").Value = 45013
$ws.Range("E159").Value = 16
$ws.Range("F159").Value = 100112021
$ws.Range("G159").Value = "Ají"
$ws.Range("H159").Value = "Cacho cabra verde"
$ws.Range("I159").Value = "Primera"
$ws.Range("J159").Value = 40
$ws.Range("K159").Value = 15000
$ws.Range("L159").Value = 15000
$ws.Range("M159").Value = 15000
$ws.Range("N159").Value = "`$/saco 25 kilos"
$ws.Range("O159").Value = "Región del Maule"
$ws.Range("P159").Value = 600
$ws.Range("Q159").Value = 25
$ws.Range("R159").Value = "Hortaliza"
